# -----------------------------------------------------------------------
# update user route and refactor code
#
#  - rename "User Info" -> "Get user"
#  - add a new "Update user" sheet (PATCH /user) after it, becomes active
#  - tidy up the stored selections / scroll position on a few sheets
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the existing "User Info" sheet to "Get user"
# ---------------------------------------------------------------------
$wsGet = $wb.Worksheets.Item("User Info")
$wsGet.Name = "Get user"

# ---------------------------------------------------------------------
# 2. Duplicate it to create the new "Update user" sheet right after it
#    (this also copies over all of the formatting/column widths/styles)
# ---------------------------------------------------------------------
$wsGet.Copy($null, $wsGet)
$wsUpd = $wb.Worksheets.Item($wsGet.Index + 1)
$wsUpd.Name = "Update user"

# ---------------------------------------------------------------------
# 3. Fix up the "Update user" sheet contents
# ---------------------------------------------------------------------

# route header
$wsUpd.Range("A2").Value = 'PATCH("/user")'

# clear out the copied "GET" success row (old F6:G6) - it gets rebuilt below
$wsUpd.Range("F6:G6").Clear()
$wsUpd.Rows.Item(6).EntireRow.AutoFit()

# request-body parameter rows, copied (with formatting) from "register"
$wsRegister = $wb.Worksheets.Item("register")
$wsRegister.Range("B5:E5").Copy($wsUpd.Range("B5:E5"))
$wsRegister.Range("B6:E6").Copy($wsUpd.Range("B6:E6"))

# Success row (F9:G9) - reuse the label style from the old GET row,
# then set the new response body text
$wsGet.Range("F6").Copy($wsUpd.Range("F9"))
$wsGet.Range("G6").Copy($wsUpd.Range("G9"))
$wsUpd.Rows.Item(9).RowHeight = 138
$wsUpd.Range("G9").Value = @"
{
    "status": 200,
    "data": {
        "id": "614155a27c7027b81a84f74d",
        "email": "ntikhoa321@gmail.com",
        "username": "NTIKHOA"
    },
    "error": null,
    "message": "Update user successfully"
}
"@

# Failure row (F10:G10) - same "email already exists" failure as register's
$wsRegister.Range("F12:G12").Copy($wsUpd.Range("F10:G10"))
$wsUpd.Rows.Item(10).RowHeight = 82.8

# column widths for E/F/G on the new sheet (target stored widths are
# 12.59765625 / 15.69921875 / 36.3984375 - the engine's ColumnWidth setter
# re-adds ~5/7 of padding on save, so back that out here)
$wsUpd.Columns.Item(5).ColumnWidth = 11.857142857142858
$wsUpd.Columns.Item(6).ColumnWidth = 15.0
$wsUpd.Columns.Item(7).ColumnWidth = 35.714285714285715

# keep the "portrait" page setup (like the other API-route sheets have)
$wsUpd.PageSetup.Orientation = 1

# sheet view: no frozen/top-left scroll, land on F9
$wsUpd.Range("F9").Select()

# ---------------------------------------------------------------------
# 4. Touch up view state (scroll/selection) left behind on other sheets
# ---------------------------------------------------------------------

# register: scrolled down a bit, F12:G12 selected
$wsRegister.Application.Goto($wsRegister.Range("A3"))
$wsRegister.Range("F12:G12").Select()

# login: selection collapsed to a single cell
$wsLogin = $wb.Worksheets.Item("login")
$wsLogin.Range("F8").Select()

# Get user: no longer frozen at the top, collapsed selection moved down
$wsGet.Range("F7:G7").Select()

# ---------------------------------------------------------------------
# 5. Make "Update user" the active sheet/tab
# ---------------------------------------------------------------------
$wsUpd.Select()
$wsUpd.Range("F9").Select()
